$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on column D for this sheet's data rows so that
# numeric-looking price strings (e.g. "1.00", "316.63") are stored as literal
# text, matching the source file's inlineStr cells instead of being coerced
# to Excel numbers (which would drop formatting like trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.743.12"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.466.28"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.63"
$ws.Range("E5").Value = "  +1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.80"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  +3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.73"
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +7.39%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.849.60"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.81"
$ws.Range("E15").Value = "  +2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.465.49"
$ws.Range("E16").Value = "  -3.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.731.05"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  +2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  +2.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.66"
$ws.Range("E21").Value = "  +3.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.15"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.20"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.74"
$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  +1.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.80"
$ws.Range("E30").Value = "  -1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.01"
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0764"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("E35").Value = "  +2.51%  "

$ws.Range("E36").Value = "  -3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.103"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  -2.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -2.08%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.974.51"
$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.90"
$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.03"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.703.34"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.03"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.12"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.01"
$ws.Range("E51").Value = "  -0.38%  "
